# Update UnitPriceUSD (col I) and UnitPriceEURO (col J) values for rows 2-8
# to reflect refreshed currency conversion rates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = 44.38;  J = 47 },
    @{ Row = 3; I = 1096.11; J = 1161 },
    @{ Row = 4; I = 297.33; J = 314.93 },
    @{ Row = 5; I = 133.13; J = 141.01 },
    @{ Row = 6; I = 2387.48; J = 2528.82 },
    @{ Row = 7; I = 621.28; J = 658.0599999999999 },
    @{ Row = 8; I = 1198.18; J = 1269.11 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I   # Column I = UnitPriceUSD
    $ws.Cells.Item($u.Row, 10).Value = $u.J  # Column J = UnitPriceEURO
}
